$wb = $excel.ActiveWorkbook

# Sheet 1: 26 changes
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 352
$ws.Range("F3").Value = 1255
$ws.Range("F4").Value = 628
$ws.Range("F8").Value = 2042
$ws.Range("F9").Value = 92
$ws.Range("F10").Value = 744
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 130
$ws.Range("F16").Value = 781
$ws.Range("F17").Value = 19
$ws.Range("F18").Value = 625
$ws.Range("F19").Value = 1203
$ws.Range("F22").Value = 723
$ws.Range("F23").Value = 693
$ws.Range("F24").Value = 72
$ws.Range("F26").Value = 615
$ws.Range("F27").Value = 1180
$ws.Range("F28").Value = 104
$ws.Range("F30").Value = 4744
$ws.Range("F31").Value = 217
$ws.Range("F32").Value = 1366
$ws.Range("F33").Value = 5727
$ws.Range("F34").Value = 942
$ws.Range("F35").Value = 565
$ws.Range("F41").Value = 628
$ws.Range("F46").Value = 23

# Sheet 2: 15 changes
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 28
$ws.Range("F6").Value = 2064
$ws.Range("F9").Value = 104
$ws.Range("F10").Value = 459
$ws.Range("F12").Value = 91
$ws.Range("F13").Value = 91
$ws.Range("F14").Value = 114
$ws.Range("F16").Value = 629
$ws.Range("F17").Value = 629
$ws.Range("F19").Value = 781
$ws.Range("F22").Value = 34
$ws.Range("F29").Value = 1704
$ws.Range("F38").Value = 55
$ws.Range("F42").Value = 879
$ws.Range("F43").Value = 465

# Sheet 3: 5 changes
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 75
$ws.Range("F5").Value = 647
$ws.Range("F6").Value = 722
$ws.Range("F7").Value = 337
$ws.Range("F8").Value = 190

# Sheet 4: 35 changes
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 352
$ws.Range("F4").Value = 647
$ws.Range("F5").Value = 28
$ws.Range("F6").Value = 1255
$ws.Range("F7").Value = 722
$ws.Range("F11").Value = 337
$ws.Range("F12").Value = 190
$ws.Range("F13").Value = 190
$ws.Range("F14").Value = 2064
$ws.Range("F15").Value = 2042
$ws.Range("F17").Value = 744
$ws.Range("F18").Value = 104
$ws.Range("F20").Value = 103
$ws.Range("F21").Value = 130
$ws.Range("F23").Value = 781
$ws.Range("F24").Value = 19
$ws.Range("F25").Value = 1203
$ws.Range("F26").Value = 91
$ws.Range("F28").Value = 723
$ws.Range("F29").Value = 114
$ws.Range("F30").Value = 693
$ws.Range("F31").Value = 629
$ws.Range("F32").Value = 615
$ws.Range("F33").Value = 104
$ws.Range("F37").Value = 34
$ws.Range("F38").Value = 4745
$ws.Range("F39").Value = 1366
$ws.Range("F40").Value = 5727
$ws.Range("F41").Value = 942
$ws.Range("F42").Value = 1704
$ws.Range("F43").Value = 565
$ws.Range("F46").Value = 628
$ws.Range("F47").Value = 55
$ws.Range("F50").Value = 879
$ws.Range("F51").Value = 465
